# Commit: "improve scene script subscribe scene events"
#
# The effect_script sheet's "script_param" column is split into three
# distinct, differently-typed parameters (str_param / json_param /
# int_param) so scripts subscribing to scene events can receive a
# string, a json blob and an int instead of a single generic param.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effect_script")

# --- header rows -----------------------------------------------------
# row 10: category label repeated across the (now 3) param columns
$ws.Range("D10").Value = "参数"
$ws.Range("E10").Value = "参数"

# row 11: field names - script_param -> str_param, plus two new fields
$ws.Range("C11").Value = "str_param"
$ws.Range("D11").Value = "json_param"
$ws.Range("E11").Value = "int_param"

# row 12: field types
$ws.Range("D12").Value = "string"
$ws.Range("E12").Value = "int"

# --- data rows 13-17 ---------------------------------------------------
$ws.Range("C13").Value = "str_1"
$ws.Range("D13").Value = "{}"
$ws.Range("E13").Value = 1

$ws.Range("C14").Value = "str_2"
$ws.Range("D14").Value = "{}"
$ws.Range("E14").Value = 2

$ws.Range("C15").Value = "str_3"
$ws.Range("D15").Value = "{}"
$ws.Range("E15").Value = 3

$ws.Range("C16").Value = "str_4"
$ws.Range("D16").Value = "{}"
$ws.Range("E16").Value = 4

$ws.Range("C17").Value = "str_5"
$ws.Range("D17").Value = "{}"
$ws.Range("E17").Value = 5

# --- active sheet / selection ------------------------------------------
# Author was reviewing the effect_script sheet (tab + selection) when the
# workbook was saved, rather than effect_searcher.
$ws.Activate()
$ws.Range("D19").Select()
